$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = 24.05
$ws.Range("B4").Value = "J'ai passé la matinée sur un problème passablement frustrant. En essayant de modifier le chargeur d'image pour le faire créer des objets, j'ai oublié de convertir Image.fromarray(image) en ImageTk.PhotoImage. Problème résolu, mais du temps a été perdu inutilement sur un problème facile. Cependant, avec la nouvelle architecture orientée objet, il me sera plus facile de récupérer une partie de l'image."

$ws.Range("B5").Select()
